# Auto-generated edit script: update crypto price/volume table per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '27.048.58'
$ws.Range('D2').Style = $ws.Range('C2').Style
$ws.Range('E2').Value = '  -3.12%  '
$ws.Range('D3').Value = "'" + '1.708.27'
$ws.Range('D3').Style = $ws.Range('C3').Style
$ws.Range('E3').Value = '  -3.51%  '
$ws.Range('D4').Value = "'" + '1.003'
$ws.Range('D4').Style = $ws.Range('C4').Style
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'" + '307.91'
$ws.Range('D5').Style = $ws.Range('C5').Style
$ws.Range('E5').Value = '  -6.20%  '
$ws.Range('D6').Value = "'" + '1.003'
$ws.Range('D6').Style = $ws.Range('C6').Style
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').Value = "'" + '0.4704'
$ws.Range('D7').Style = $ws.Range('C7').Style
$ws.Range('E7').Value = '  +5.14%  '
$ws.Range('D8').Value = "'" + '0.3416'
$ws.Range('D8').Style = $ws.Range('C8').Style
$ws.Range('E8').Value = '  -3.95%  '
$ws.Range('D9').Value = "'" + '41.91'
$ws.Range('D9').Style = $ws.Range('C9').Style
$ws.Range('E9').Value = '  -0.48%  '
$ws.Range('D10').Value = "'" + '0.07248'
$ws.Range('D10').Style = $ws.Range('C10').Style
$ws.Range('E10').Value = '  -2.54%  '
$ws.Range('D11').Value = "'" + '1.033'
$ws.Range('D11').Style = $ws.Range('C11').Style
$ws.Range('E11').Value = '  -6.38%  '
$ws.Range('D12').Value = "'" + '1.002'
$ws.Range('D12').Style = $ws.Range('C12').Style
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = "'" + '19.75'
$ws.Range('D13').Style = $ws.Range('C13').Style
$ws.Range('E13').Value = '  -5.98%  '
$ws.Range('D14').Value = "'" + '5.836'
$ws.Range('D14').Style = $ws.Range('C14').Style
$ws.Range('E14').Value = '  -3.24%  '
$ws.Range('D15').Value = "'" + '1.710.29'
$ws.Range('D15').Style = $ws.Range('C15').Style
$ws.Range('E15').Value = '  -3.28%  '
$ws.Range('D16').Value = "'" + '6.825'
$ws.Range('D16').Style = $ws.Range('C16').Style
$ws.Range('E16').Value = '  -5.84%  '
$ws.Range('D17').Value = "'" + '88.96'
$ws.Range('D17').Style = $ws.Range('C17').Style
$ws.Range('E17').Value = '  -4.63%  '
$ws.Range('D18').Value = "'" + '0.00001036'
$ws.Range('D18').Style = $ws.Range('C18').Style
$ws.Range('E18').Value = '  -2.58%  '
$ws.Range('D19').Value = "'" + '0.06349'
$ws.Range('D19').Style = $ws.Range('C19').Style
$ws.Range('E19').Value = '  -1.42%  '
$ws.Range('D20').Value = "'" + '1.003'
$ws.Range('D20').Style = $ws.Range('C20').Style
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('D21').Value = "'" + '16.43'
$ws.Range('D21').Style = $ws.Range('C21').Style
$ws.Range('E21').Value = '  -4.04%  '
$ws.Range('D22').Value = "'" + '5.595'
$ws.Range('D22').Style = $ws.Range('C22').Style
$ws.Range('E22').Value = '  -3.28%  '
$ws.Range('D23').Value = "'" + '27.071.46'
$ws.Range('D23').Style = $ws.Range('C23').Style
$ws.Range('E23').Value = '  -3.18%  '
$ws.Range('D24').Value = "'" + '10.83'
$ws.Range('D24').Style = $ws.Range('C24').Style
$ws.Range('E24').Value = '  -4.10%  '
$ws.Range('D25').Value = "'" + '2.114'
$ws.Range('D25').Style = $ws.Range('C25').Style
$ws.Range('E25').Value = '  +0.29%  '
$ws.Range('D26').Value = "'" + '157.05'
$ws.Range('D26').Style = $ws.Range('C26').Style
$ws.Range('E26').Value = '  -2.99%  '
$ws.Range('D27').Value = "'" + '19.45'
$ws.Range('D27').Style = $ws.Range('C27').Style
$ws.Range('E27').Value = '  -4.70%  '
$ws.Range('D28').Value = "'" + '1.906.16'
$ws.Range('D28').Style = $ws.Range('C28').Style
$ws.Range('E28').Value = '  -3.40%  '
$ws.Range('D29').Value = "'" + '2.073'
$ws.Range('D29').Style = $ws.Range('C29').Style
$ws.Range('E29').Value = '  -4.50%  '
$ws.Range('D30').Value = "'" + '119.11'
$ws.Range('D30').Style = $ws.Range('C30').Style
$ws.Range('E30').Value = '  -4.75%  '
$ws.Range('D31').Value = "'" + '1.009'
$ws.Range('D31').Style = $ws.Range('C31').Style
$ws.Range('E31').Value = '  -9.24%  '
$ws.Range('D32').Value = "'" + '0.09082'
$ws.Range('D32').Style = $ws.Range('C32').Style
$ws.Range('E32').Value = '  -1.40%  '
$ws.Range('E33').Value = '  -2.40%  '
$ws.Range('D34').Value = "'" + '5.292'
$ws.Range('D34').Style = $ws.Range('C34').Style
$ws.Range('E34').Value = '  -6.42%  '
$ws.Range('D35').Value = "'" + '0.02189'
$ws.Range('D35').Style = $ws.Range('C35').Style
$ws.Range('E35').Value = '  -4.49%  '
$ws.Range('E36').Value = '  -5.77%  '
$ws.Range('D37').Value = "'" + '10.98'
$ws.Range('D37').Style = $ws.Range('C37').Style
$ws.Range('E37').Value = '  -7.56%  '
$ws.Range('E38').Value = '  -5.62%  '
$ws.Range('B39').Value = 'Frax'
$ws.Range('C39').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D39').Value = "'" + '1.002'
$ws.Range('D39').Style = $ws.Range('C39').Style
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').Value = "'" + '4.722'
$ws.Range('D40').Style = $ws.Range('C40').Style
$ws.Range('E40').Value = '  -5.08%  '
$ws.Range('B41').Value = 'WEMIXTOKEN'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').Value = "'" + '1.397'
$ws.Range('D41').Style = $ws.Range('C41').Style
$ws.Range('E41').Value = '  +0.30%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = "'" + '0.5855'
$ws.Range('D42').Style = $ws.Range('C42').Style
$ws.Range('E42').Value = '  -7.53%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = "'" + '1.095'
$ws.Range('D43').Style = $ws.Range('C43').Style
$ws.Range('E43').Value = '  -7.53%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = "'" + '7.454'
$ws.Range('D44').Style = $ws.Range('C44').Style
$ws.Range('E44').Value = '  -5.85%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = "'" + '12.56'
$ws.Range('D45').Style = $ws.Range('C45').Style
$ws.Range('E45').Value = '  -5.39%  '
$ws.Range('D46').Value = "'" + '0.5636'
$ws.Range('D46').Style = $ws.Range('C46').Style
$ws.Range('E46').Value = '  -4.35%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').Value = "'" + '3.553'
$ws.Range('D47').Style = $ws.Range('C47').Style
$ws.Range('E47').Value = '  -5.17%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = "'" + '116.82'
$ws.Range('D48').Style = $ws.Range('C48').Style
$ws.Range('E48').Value = '  -4.90%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = "'" + '1.831'
$ws.Range('D49').Style = $ws.Range('C49').Style
$ws.Range('E49').Value = '  -6.70%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'" + '0.06617'
$ws.Range('D50').Style = $ws.Range('C50').Style
$ws.Range('E50').Value = '  -4.21%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').Value = "'" + '1.080'
$ws.Range('D51').Style = $ws.Range('C51').Style
$ws.Range('E51').Value = '  -5.28%  '
